$wb = $excel.ActiveWorkbook

# --- 1. Update the ColumnHeaders sheet's definition for iode_quality_flag ---
# (was the terse "flag", now a fuller description)
$colHeaders = $wb.Worksheets.Item("ColumnHeaders")
$colHeaders.Range("B13").Value = "IODE Quality Flag primary level"

# --- 2. Insert new "CategoricalVariables" sheet before ColumnHeaders ---
$catVars = $wb.Worksheets.Add()
$catVars.Name = "CategoricalVariables"

# Re-acquire the ColumnHeaders reference: the sheet collection shifted when
# the new sheet was inserted ahead of it, so the old handle is stale.
$colHeaders = $wb.Worksheets.Item("ColumnHeaders")

# --- 3. Populate CategoricalVariables data grid (A1:C8) ---
$catVars.Range("A1").Value = "attributeName"
$catVars.Range("B1").Value = "code"
$catVars.Range("C1").Value = "definition"

$catVars.Range("A2").Value = "toi_source"
$catVars.Range("B2").Value = "toi_niskin"
$catVars.Range("C2").Value = "sample bottle was filled from a Niskin bottle on CTD rosette"

$catVars.Range("A3").Value = "toi_source"
$catVars.Range("B3").Value = "toi_underway"
$catVars.Range("C3").Value = "sample bottle was filled from the ship's underway system"

$catVars.Range("A4").Value = "iode_quality_flag"
$catVars.Range("B4").Value = 1
$catVars.Range("C4").Value = "good"

$catVars.Range("A5").Value = "iode_quality_flag"
$catVars.Range("B5").Value = 2
$catVars.Range("C5").Value = "quality not evaluated, not available or unknown"

$catVars.Range("A6").Value = "iode_quality_flag"
$catVars.Range("B6").Value = 3
$catVars.Range("C6").Value = "questionable/suspect"

$catVars.Range("A7").Value = "iode_quality_flag"
$catVars.Range("B7").Value = 4
$catVars.Range("C7").Value = "bad"

$catVars.Range("A8").Value = "iode_quality_flag"
$catVars.Range("B8").Value = 9
$catVars.Range("C8").Value = "missing data"

# --- 4. Restore the UI selection state on ColumnHeaders (no longer the active tab) ---
[void]$colHeaders.Range("B15").Select()

# --- 5. Make CategoricalVariables the active tab again, with its own lingering selection ---
$catVars = $wb.Worksheets.Item("CategoricalVariables")
[void]$catVars.Range("B32").Select()
